$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Förändrad) becomes 46066 for every data row (2-14)
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 3).Value = 46066
}

# Rows 3-14 get their A (Beteckning), B (Datum), G (Area) values reassigned
# per the permutation described in the diff.
$ws.Cells.Item(3, 1).Value = "A 1053-2022"
$ws.Cells.Item(3, 2).Value = 44571
$ws.Cells.Item(3, 7).Value = 1.7

$ws.Cells.Item(4, 1).Value = "A 24221-2021"
$ws.Cells.Item(4, 2).Value = 44336.78922453704
$ws.Cells.Item(4, 7).Value = 1.1

$ws.Cells.Item(5, 1).Value = "A 46826-2022"
$ws.Cells.Item(5, 2).Value = 44848
$ws.Cells.Item(5, 7).Value = 4.5

$ws.Cells.Item(6, 1).Value = "A 50934-2024"
$ws.Cells.Item(6, 2).Value = 45602
$ws.Cells.Item(6, 7).Value = 0.6

$ws.Cells.Item(7, 1).Value = "A 31120-2023"
$ws.Cells.Item(7, 2).Value = 45113
$ws.Cells.Item(7, 7).Value = 0.2

$ws.Cells.Item(8, 1).Value = "A 46779-2025"
$ws.Cells.Item(8, 2).Value = 45926
$ws.Cells.Item(8, 7).Value = 1.5

$ws.Cells.Item(9, 1).Value = "A 64431-2023"
$ws.Cells.Item(9, 2).Value = 45280
$ws.Cells.Item(9, 7).Value = 0.5

$ws.Cells.Item(10, 1).Value = "A 56948-2025"
$ws.Cells.Item(10, 2).Value = 45978.64356481482
$ws.Cells.Item(10, 7).Value = 4.7

$ws.Cells.Item(11, 1).Value = "A 56917-2025"
$ws.Cells.Item(11, 2).Value = 45978.58453703704
$ws.Cells.Item(11, 7).Value = 0.7

$ws.Cells.Item(12, 1).Value = "A 27724-2022"
$ws.Cells.Item(12, 2).Value = 44743.48386574074
$ws.Cells.Item(12, 7).Value = 1.3

$ws.Cells.Item(13, 1).Value = "A 8748-2022"
$ws.Cells.Item(13, 2).Value = 44613
$ws.Cells.Item(13, 7).Value = 1

$ws.Cells.Item(14, 1).Value = "A 64445-2023"
$ws.Cells.Item(14, 2).Value = 45280
$ws.Cells.Item(14, 7).Value = 3.7
